$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Window position (best effort; cosmetic window-on-screen coordinates) ---
try {
    $excel.ActiveWindow.Left = 3120
    $excel.ActiveWindow.Top = 3120
} catch {
}

# --- Row 4: Mini Factory (迷你工厂) ---
# Effect text changes to the new "graveyard" version (becomes new shared string, taking slot vacated by old E4 text).
$ws.Range("E4").Value = "回合结束时：本牌点数减1。然后将墓地第1张“机器人”牌放在本牌所在槽位中，本牌在备战区时，可以改为选墓地1张“机器人”牌放在本牌所在槽位中。本牌点数因此降低至0时，本牌不会死亡而是弃置。<br>`n从手牌发动：将本牌放到房间区任意位置。"
# Row height grows from 42.75 to 57 to fit the longer wrapped text.
$ws.Rows(4).RowHeight = 57

# --- Row 5: Kamikaze bot (自爆机器人) ---
# Effect text (E5) is left untouched -- its shared-string index simply shifts down once the
# orphaned old-E4 string is compacted out of the table.
# Cost/rank changes from 2 to 1.
$ws.Range("C5").Value = 1

# --- Row 6: Watchdog/Sentinel bot (哨戒机器人) ---
# Effect text changes to the new combined "discard pile + move to adjacent" version.
$ws.Range("E6").Value = "回合结束时：本牌点数减1，然后将弃牌堆第1张“机器人”牌放在本牌所在槽位中，本牌在备战区时，可以改为选弃牌堆1张“机器人”牌放在本牌所在槽位中。这之后，将本牌移动到1个相邻槽位中。本牌点数因此降低至0时，本牌不会死亡而是弃置。<br>`n从手牌发动：将本牌放到房间区任意位置。"
# Cost/rank changes from 1 to 2.
$ws.Range("C6").Value = 2

# --- Selection moves from D5 to E5 ---
$ws.Range("E5").Select()
